$d = $word.ActiveDocument

# Locate the last paragraph (the one right after "Opdracht 4:"), which currently
# holds a single run containing only a tab, followed by the (hidden) _GoBack bookmark.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$tabRange = $r.Duplicate
$tabRange.MoveEnd(1, -1)

$tabStart = $tabRange.Start

# Remove the existing lone-tab run; the (position-anchored) bookmark stays put.
$tabRange.Delete()

# Re-insert the tab together with the new explanatory runs (incl. the spell-check
# markers around the misspelled "datatpe") right where the tab used to be, so the
# bookmark ends up after all of the new content, matching the original authoring flow.
$insPoint = $d.Range($tabStart, $tabStart)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:r><w:t>Het verschil waarom de verwijzing van j naar i bij het primitieve type wel gelijk blijft maar voor het object niet zit hem in het feit dat je met een primitief type feitelijk een kopie maakt van de parameter</w:t></w:r><w:r><w:t xml:space="preserve"> (dit is een directe opslag van de waarde in het ‘vakje’ van de variabele met het primitieve </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>datatpe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">, terwijl je bij het object met zowel a als b een verwijzing naar hetzelfde object maakt. </w:t></w:r><w:r><w:t>De waarde van het object wordt dus niet rechtstreeks in de variabele opgeslagen, dit is slechts een referentie naar het object.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($xml)
